$wb = $excel.ActiveWorkbook

# ---- Sheet: GLOBAL RESULTS ----
$ws = $wb.Worksheets.Item('GLOBAL RESULTS')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg structure MAC'
$ws.Cells.Item(2, 2).Value = '%'
$ws.Cells.Item(2, 3).Value = 54.44126671372609
$ws.Cells.Item(3, 1).Value = 'Xcg structure BRF'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 12.437815810838156
$ws.Cells.Item(4, 1).Value = 'Zcg structure MAC'
$ws.Cells.Item(4, 2).Value = '%'
$ws.Cells.Item(4, 3).Value = 27.67850347037711
$ws.Cells.Item(5, 1).Value = 'Zcg structure BRF'
$ws.Cells.Item(5, 2).Value = 'm'
$ws.Cells.Item(5, 3).Value = 0.637274225272874
$ws.Cells.Item(6, 1).Value = ' '
$ws.Cells.Item(6, 2).Value = $null
$ws.Cells.Item(6, 3).Value = $null
$ws.Cells.Item(7, 1).Value = 'Xcg structure and engines MAC'
$ws.Cells.Item(7, 2).Value = '%'
$ws.Cells.Item(7, 3).Value = 32.206395870082275
$ws.Cells.Item(8, 1).Value = 'Xcg structure and engines BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 11.925876640929975
$ws.Cells.Item(9, 1).Value = 'Zcg structure and engines MAC'
$ws.Cells.Item(9, 2).Value = '%'
$ws.Cells.Item(9, 3).Value = 30.784599937879275
$ws.Cells.Item(10, 1).Value = 'Zcg structure and engines BRF'
$ws.Cells.Item(10, 2).Value = 'm'
$ws.Cells.Item(10, 3).Value = 0.708789479776021
$ws.Cells.Item(11, 1).Value = ' '
$ws.Cells.Item(11, 2).Value = $null
$ws.Cells.Item(11, 3).Value = $null
$ws.Cells.Item(12, 1).Value = 'Xcg operating empty mass MAC'
$ws.Cells.Item(12, 2).Value = '%'
$ws.Cells.Item(12, 3).Value = 32.206395870082275
$ws.Cells.Item(13, 1).Value = 'Xcg operating empty mass BRF'
$ws.Cells.Item(13, 2).Value = 'm'
$ws.Cells.Item(13, 3).Value = 11.925876640929975
$ws.Cells.Item(14, 1).Value = 'Zcg operating empty mass MAC'
$ws.Cells.Item(14, 2).Value = '%'
$ws.Cells.Item(14, 3).Value = 30.784599937879275
$ws.Cells.Item(15, 1).Value = 'Zcg operating empty mass BRF'
$ws.Cells.Item(15, 2).Value = 'm'
$ws.Cells.Item(15, 3).Value = 0.708789479776021
$ws.Cells.Item(16, 1).Value = ' '
$ws.Cells.Item(16, 2).Value = $null
$ws.Cells.Item(16, 3).Value = $null
$ws.Cells.Item(17, 1).Value = 'Xcg maximum zero fuel mass MAC'
$ws.Cells.Item(17, 2).Value = '%'
$ws.Cells.Item(17, 3).Value = 34.28982845054642
$ws.Cells.Item(18, 1).Value = 'Xcg maximum zero fuel mass BRF'
$ws.Cells.Item(18, 2).Value = 'm'
$ws.Cells.Item(18, 3).Value = 11.973845920954258
$ws.Cells.Item(19, 1).Value = 'Zcg maximum zero fuel mass MAC'
$ws.Cells.Item(19, 2).Value = '%'
$ws.Cells.Item(19, 3).Value = 20.6733270276342
$ws.Cells.Item(20, 1).Value = 'Zcg maximum zero fuel mass BRF'
$ws.Cells.Item(20, 2).Value = 'm'
$ws.Cells.Item(20, 3).Value = 0.4759859390320157
$ws.Cells.Item(21, 1).Value = ' '
$ws.Cells.Item(21, 2).Value = $null
$ws.Cells.Item(21, 3).Value = $null
$ws.Cells.Item(22, 1).Value = 'Xcg maximum take-off mass MAC'
$ws.Cells.Item(22, 2).Value = '%'
$ws.Cells.Item(22, 3).Value = 54.284755273589916
$ws.Cells.Item(23, 1).Value = 'Xcg maximum take-off mass BRF'
$ws.Cells.Item(23, 2).Value = 'm'
$ws.Cells.Item(23, 3).Value = 12.434212266778562
$ws.Cells.Item(24, 1).Value = 'Zcg maximum take-off mass MAC'
$ws.Cells.Item(24, 2).Value = '%'
$ws.Cells.Item(24, 3).Value = 27.842825670894356
$ws.Cells.Item(25, 1).Value = 'Zcg maximum take-off mass BRF'
$ws.Cells.Item(25, 2).Value = 'm'
$ws.Cells.Item(25, 3).Value = 0.6410576055102428
$ws.Cells.Item(26, 1).Value = ' '
$ws.Cells.Item(26, 2).Value = $null
$ws.Cells.Item(26, 3).Value = $null
$ws.Cells.Item(27, 1).Value = 'Max forward Xcg MAC'
$ws.Cells.Item(27, 2).Value = '%'
$ws.Cells.Item(27, 3).Value = 18.282686812078495
$ws.Cells.Item(28, 1).Value = 'Max aft Xcg MAC'
$ws.Cells.Item(28, 2).Value = '%'
$ws.Cells.Item(28, 3).Value = 54.284755273589916

# ---- Sheet: FUSELAGE ----
$ws = $wb.Worksheets.Item('FUSELAGE')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg LRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 11.04753745546374
$ws.Cells.Item(3, 1).Value = 'Ycg LRF'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(4, 1).Value = 'Zcg LRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(6, 1).Value = 'Xcg BRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 11.04753745546374
$ws.Cells.Item(7, 1).Value = 'Ycg BRF'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 0.0
$ws.Cells.Item(8, 1).Value = 'Zcg BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 0.0
$ws.Cells.Item(9, 1).Value = ' '
$ws.Cells.Item(9, 2).Value = $null
$ws.Cells.Item(9, 3).Value = $null
$ws.Cells.Item(10, 1).Value = 'Xcg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(10, 2).Value = $null
$ws.Cells.Item(10, 3).Value = $null
$ws.Cells.Item(11, 1).Value = 'SFORZA'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 11.500334910927485
$ws.Cells.Item(12, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(12, 2).Value = 'm'
$ws.Cells.Item(12, 3).Value = 10.594739999999998

# ---- Sheet: WING ----
$ws = $wb.Worksheets.Item('WING')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg LRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 1.330689739372911
$ws.Cells.Item(3, 1).Value = 'Ycg LRF (semi-wing)'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 4.735499999999998
$ws.Cells.Item(4, 1).Value = 'Zcg LRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(6, 1).Value = 'Xcg BRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 12.33068973937291
$ws.Cells.Item(7, 1).Value = 'Ycg BRF (semi-wing)'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 4.735499999999997
$ws.Cells.Item(8, 1).Value = 'Zcg BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 1.5999999999999999
$ws.Cells.Item(9, 1).Value = ' '
$ws.Cells.Item(9, 2).Value = $null
$ws.Cells.Item(9, 3).Value = $null
$ws.Cells.Item(10, 1).Value = 'Xcg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(10, 2).Value = $null
$ws.Cells.Item(10, 3).Value = $null
$ws.Cells.Item(11, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 1.330689739372911
$ws.Cells.Item(12, 1).Value = ' '
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(13, 1).Value = 'Ycg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(13, 2).Value = $null
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(14, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(14, 2).Value = 'm'
$ws.Cells.Item(14, 3).Value = 4.735499999999998

# ---- Sheet: FUEL TANK ----
$ws = $wb.Worksheets.Item('FUEL TANK')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg LRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 1.0717832283332958
$ws.Cells.Item(3, 1).Value = 'Ycg LRF'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(4, 1).Value = 'Zcg LRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(6, 1).Value = 'Xcg BRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 12.071783228333294
$ws.Cells.Item(7, 1).Value = 'Ycg BRF'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 0.0
$ws.Cells.Item(8, 1).Value = 'Zcg BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 1.5999999999999999
$ws.Cells.Item(9, 1).Value = ' '
$ws.Cells.Item(9, 2).Value = $null
$ws.Cells.Item(9, 3).Value = $null

# ---- Sheet: HORIZONTAL TAIL ----
$ws = $wb.Worksheets.Item('HORIZONTAL TAIL')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg LRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 0.7599657170108016
$ws.Cells.Item(3, 1).Value = 'Ycg LRF (semi-tail)'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 1.3888239999999996
$ws.Cells.Item(4, 1).Value = 'Zcg LRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(6, 1).Value = 'Xcg BRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 26.059965717010797
$ws.Cells.Item(7, 1).Value = 'Ycg BRF (semi-tail)'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 1.3888239999999994
$ws.Cells.Item(8, 1).Value = 'Zcg BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 5.737399999999999
$ws.Cells.Item(9, 1).Value = ' '
$ws.Cells.Item(9, 2).Value = $null
$ws.Cells.Item(9, 3).Value = $null
$ws.Cells.Item(10, 1).Value = 'Xcg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(10, 2).Value = $null
$ws.Cells.Item(10, 3).Value = $null
$ws.Cells.Item(11, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 0.7599657170108016
$ws.Cells.Item(12, 1).Value = ' '
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(13, 1).Value = 'Ycg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(13, 2).Value = $null
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(14, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(14, 2).Value = 'm'
$ws.Cells.Item(14, 3).Value = 1.3888239999999996

# ---- Sheet: VERTICAL TAIL ----
$ws = $wb.Worksheets.Item('VERTICAL TAIL')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg LRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 2.5640710746536337
$ws.Cells.Item(3, 1).Value = 'Ycg LRF (semi-tail)'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 1.8316
$ws.Cells.Item(4, 1).Value = 'Zcg LRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(6, 1).Value = 'Xcg BRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 24.16407107465363
$ws.Cells.Item(7, 1).Value = 'Ycg BRF (semi-tail)'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 0.0
$ws.Cells.Item(8, 1).Value = 'Zcg BRF'
$ws.Cells.Item(8, 2).Value = 'm'
$ws.Cells.Item(8, 3).Value = 3.1315999999999997
$ws.Cells.Item(9, 1).Value = ' '
$ws.Cells.Item(9, 2).Value = $null
$ws.Cells.Item(9, 3).Value = $null
$ws.Cells.Item(10, 1).Value = 'Xcg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(10, 2).Value = $null
$ws.Cells.Item(10, 3).Value = $null
$ws.Cells.Item(11, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 2.5640710746536337
$ws.Cells.Item(12, 1).Value = ' '
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(13, 1).Value = 'Ycg ESTIMATION METHOD COMPARISON'
$ws.Cells.Item(13, 2).Value = $null
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(14, 1).Value = 'TORENBEEK_1982'
$ws.Cells.Item(14, 2).Value = 'm'
$ws.Cells.Item(14, 3).Value = 1.8316

# ---- Sheet: NACELLES ----
$ws = $wb.Worksheets.Item('NACELLES')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'BALANCE ESTIMATION FOR EACH NACELLE'
$ws.Cells.Item(2, 2).Value = $null
$ws.Cells.Item(2, 3).Value = $null
$ws.Cells.Item(3, 1).Value = ' '
$ws.Cells.Item(3, 2).Value = $null
$ws.Cells.Item(3, 3).Value = $null
$ws.Cells.Item(4, 1).Value = 'NACELLE 1'
$ws.Cells.Item(4, 2).Value = $null
$ws.Cells.Item(4, 3).Value = $null
$ws.Cells.Item(5, 1).Value = 'Xcg LRF'
$ws.Cells.Item(5, 2).Value = 'm'
$ws.Cells.Item(5, 3).Value = 1.7483999999999997
$ws.Cells.Item(6, 1).Value = 'Ycg LRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 0.0
$ws.Cells.Item(7, 1).Value = 'Zcg LRF'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 0.0
$ws.Cells.Item(8, 1).Value = ' '
$ws.Cells.Item(8, 2).Value = $null
$ws.Cells.Item(8, 3).Value = $null
$ws.Cells.Item(9, 1).Value = 'Xcg BRF'
$ws.Cells.Item(9, 2).Value = 'm'
$ws.Cells.Item(9, 3).Value = 10.317419999999998
$ws.Cells.Item(10, 1).Value = 'Ycg BRF'
$ws.Cells.Item(10, 2).Value = 'm'
$ws.Cells.Item(10, 3).Value = 4.573799999999999
$ws.Cells.Item(11, 1).Value = 'Zcg BRF'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 1.0289499999999998
$ws.Cells.Item(12, 1).Value = ' '
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(13, 1).Value = ' '
$ws.Cells.Item(13, 2).Value = $null
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(14, 1).Value = ' '
$ws.Cells.Item(14, 2).Value = $null
$ws.Cells.Item(14, 3).Value = $null
$ws.Cells.Item(15, 1).Value = 'NACELLE 2'
$ws.Cells.Item(15, 2).Value = $null
$ws.Cells.Item(15, 3).Value = $null
$ws.Cells.Item(16, 1).Value = 'Xcg LRF'
$ws.Cells.Item(16, 2).Value = 'm'
$ws.Cells.Item(16, 3).Value = 1.7483999999999997
$ws.Cells.Item(17, 1).Value = 'Ycg LRF'
$ws.Cells.Item(17, 2).Value = 'm'
$ws.Cells.Item(17, 3).Value = 0.0
$ws.Cells.Item(18, 1).Value = 'Zcg LRF'
$ws.Cells.Item(18, 2).Value = 'm'
$ws.Cells.Item(18, 3).Value = 0.0
$ws.Cells.Item(19, 1).Value = ' '
$ws.Cells.Item(19, 2).Value = $null
$ws.Cells.Item(19, 3).Value = $null
$ws.Cells.Item(20, 1).Value = 'Xcg BRF'
$ws.Cells.Item(20, 2).Value = 'm'
$ws.Cells.Item(20, 3).Value = 10.317419999999998
$ws.Cells.Item(21, 1).Value = 'Ycg BRF'
$ws.Cells.Item(21, 2).Value = 'm'
$ws.Cells.Item(21, 3).Value = -4.573799999999999
$ws.Cells.Item(22, 1).Value = 'Zcg BRF'
$ws.Cells.Item(22, 2).Value = 'm'
$ws.Cells.Item(22, 3).Value = 1.0289499999999998
$ws.Cells.Item(23, 1).Value = ' '
$ws.Cells.Item(23, 2).Value = $null
$ws.Cells.Item(23, 3).Value = $null
$ws.Cells.Item(24, 1).Value = ' '
$ws.Cells.Item(24, 2).Value = $null
$ws.Cells.Item(24, 3).Value = $null
$ws.Cells.Item(25, 1).Value = ' '
$ws.Cells.Item(25, 2).Value = $null
$ws.Cells.Item(25, 3).Value = $null

# ---- Sheet: POWER PLANT ----
$ws = $wb.Worksheets.Item('POWER PLANT')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'BALANCE ESTIMATION FOR EACH ENGINE'
$ws.Cells.Item(2, 2).Value = $null
$ws.Cells.Item(2, 3).Value = $null
$ws.Cells.Item(3, 1).Value = ' '
$ws.Cells.Item(3, 2).Value = $null
$ws.Cells.Item(3, 3).Value = $null
$ws.Cells.Item(4, 1).Value = 'ENGINE 1'
$ws.Cells.Item(4, 2).Value = $null
$ws.Cells.Item(4, 3).Value = $null
$ws.Cells.Item(5, 1).Value = 'Xcg LRF'
$ws.Cells.Item(5, 2).Value = 'm'
$ws.Cells.Item(5, 3).Value = 1.0649999999999995
$ws.Cells.Item(6, 1).Value = 'Ycg LRF'
$ws.Cells.Item(6, 2).Value = 'm'
$ws.Cells.Item(6, 3).Value = 0.0
$ws.Cells.Item(7, 1).Value = 'Zcg LRF'
$ws.Cells.Item(7, 2).Value = 'm'
$ws.Cells.Item(7, 3).Value = 0.0
$ws.Cells.Item(8, 1).Value = ' '
$ws.Cells.Item(8, 2).Value = $null
$ws.Cells.Item(8, 3).Value = $null
$ws.Cells.Item(9, 1).Value = 'Xcg BRF'
$ws.Cells.Item(9, 2).Value = 'm'
$ws.Cells.Item(9, 3).Value = 9.634019999999998
$ws.Cells.Item(10, 1).Value = 'Ycg BRF'
$ws.Cells.Item(10, 2).Value = 'm'
$ws.Cells.Item(10, 3).Value = 4.573799999999999
$ws.Cells.Item(11, 1).Value = 'Zcg BRF'
$ws.Cells.Item(11, 2).Value = 'm'
$ws.Cells.Item(11, 3).Value = 1.0289499999999998
$ws.Cells.Item(12, 1).Value = ' '
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(13, 1).Value = ' '
$ws.Cells.Item(13, 2).Value = $null
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(14, 1).Value = ' '
$ws.Cells.Item(14, 2).Value = $null
$ws.Cells.Item(14, 3).Value = $null
$ws.Cells.Item(15, 1).Value = 'ENGINE 2'
$ws.Cells.Item(15, 2).Value = $null
$ws.Cells.Item(15, 3).Value = $null
$ws.Cells.Item(16, 1).Value = 'Xcg LRF'
$ws.Cells.Item(16, 2).Value = 'm'
$ws.Cells.Item(16, 3).Value = 1.0649999999999995
$ws.Cells.Item(17, 1).Value = 'Ycg LRF'
$ws.Cells.Item(17, 2).Value = 'm'
$ws.Cells.Item(17, 3).Value = 0.0
$ws.Cells.Item(18, 1).Value = 'Zcg LRF'
$ws.Cells.Item(18, 2).Value = 'm'
$ws.Cells.Item(18, 3).Value = 0.0
$ws.Cells.Item(19, 1).Value = ' '
$ws.Cells.Item(19, 2).Value = $null
$ws.Cells.Item(19, 3).Value = $null
$ws.Cells.Item(20, 1).Value = 'Xcg BRF'
$ws.Cells.Item(20, 2).Value = 'm'
$ws.Cells.Item(20, 3).Value = 9.634019999999998
$ws.Cells.Item(21, 1).Value = 'Ycg BRF'
$ws.Cells.Item(21, 2).Value = 'm'
$ws.Cells.Item(21, 3).Value = -4.573799999999999
$ws.Cells.Item(22, 1).Value = 'Zcg BRF'
$ws.Cells.Item(22, 2).Value = 'm'
$ws.Cells.Item(22, 3).Value = 1.0289499999999998
$ws.Cells.Item(23, 1).Value = ' '
$ws.Cells.Item(23, 2).Value = $null
$ws.Cells.Item(23, 3).Value = $null
$ws.Cells.Item(24, 1).Value = ' '
$ws.Cells.Item(24, 2).Value = $null
$ws.Cells.Item(24, 3).Value = $null
$ws.Cells.Item(25, 1).Value = ' '
$ws.Cells.Item(25, 2).Value = $null
$ws.Cells.Item(25, 3).Value = $null

# ---- Sheet: LANDING GEARS ----
$ws = $wb.Worksheets.Item('LANDING GEARS')
$ws.Cells.Item(1, 1).Value = 'Description'
$ws.Cells.Item(1, 2).Value = 'Unit'
$ws.Cells.Item(1, 3).Value = 'Value'
$ws.Cells.Item(2, 1).Value = 'Xcg BRF'
$ws.Cells.Item(2, 2).Value = 'm'
$ws.Cells.Item(2, 3).Value = 12.321708214766055
$ws.Cells.Item(3, 1).Value = 'Ycg BRF'
$ws.Cells.Item(3, 2).Value = 'm'
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(4, 1).Value = 'Zcg BRF'
$ws.Cells.Item(4, 2).Value = 'm'
$ws.Cells.Item(4, 3).Value = -1.8746386998784927
$ws.Cells.Item(5, 1).Value = ' '
$ws.Cells.Item(5, 2).Value = $null
$ws.Cells.Item(5, 3).Value = $null
